# Add data for 2022-07-15
# Updates the "Through 2022-07-06" -> "Through 2022-07-07" sheet/title and
# shared-string header, plus new/updated carjacking counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (and update its header cell text in B1)
$ws.Name = "Through 2022-07-07"
$ws.Range("B1").Value = "July 2022 (through July 07)"

# Column B = "July 2022" (current month-to-date column)
$ws.Range("B8").Value = 2    # North Lawndale
$ws.Range("B16").Value = 3   # Washington Heights
$ws.Range("B34").Value = 1   # Riverdale (new)

# Column I = "July 2021"
$ws.Range("I2").Value = 4    # Austin
$ws.Range("I12").Value = 1   # Hyde Park (new)
$ws.Range("I44").Value = 1   # New City (new)
$ws.Range("I53").Value = 2   # Calumet Heights
$ws.Range("I71").Value = 1   # Galewood (new)
$ws.Range("I78").Value = 1   # Lake View (new)

# Column P = "July 2020"
$ws.Range("P5").Value = 1    # Garfield Park (new)
$ws.Range("P53").Value = 2   # Calumet Heights

# Column AK = "July 2017"
$ws.Range("AK3").Value = 1   # Englewood (new)
$ws.Range("AK38").Value = 1  # West Town (new)
$ws.Range("AK39").Value = 1  # Wicker Park (new)
